$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 194
$ws1.Range("F3").Value = 111
$ws1.Range("F5").Value = 971
$ws1.Range("F6").Value = 5336
$ws1.Range("F7").Value = 459
$ws1.Range("F8").Value = 646
$ws1.Range("F9").Value = 921
$ws1.Range("F10").Value = 820
$ws1.Range("F12").Value = 32
$ws1.Range("F14").Value = 21
$ws1.Range("F17").Value = 1780
$ws1.Range("F19").Value = 847
$ws1.Range("F22").Value = 315
$ws1.Range("F23").Value = 522
$ws1.Range("F25").Value = 1048
$ws1.Range("F28").Value = 2657
$ws1.Range("F29").Value = 174
$ws1.Range("F31").Value = 56
$ws1.Range("F32").Value = 93
$ws1.Range("F33").Value = 28
$ws1.Range("F34").Value = 316
$ws1.Range("F39").Value = 278
$ws1.Range("F40").Value = 652
$ws1.Range("F43").Value = 46

# Sheet "全部类型" (index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 194
$ws4.Range("F4").Value = 111
$ws4.Range("F5").Value = 971
$ws4.Range("F7").Value = 5336
$ws4.Range("F8").Value = 459
$ws4.Range("F9").Value = 646
$ws4.Range("F12").Value = 921
$ws4.Range("F13").Value = 820
$ws4.Range("F17").Value = 32
$ws4.Range("F19").Value = 21
$ws4.Range("F23").Value = 1780
$ws4.Range("F25").Value = 847
$ws4.Range("F27").Value = 315
$ws4.Range("F29").Value = 522
$ws4.Range("F31").Value = 1048
$ws4.Range("F32").Value = 2657
$ws4.Range("F33").Value = 174
$ws4.Range("F35").Value = 56
$ws4.Range("F36").Value = 93
$ws4.Range("F37").Value = 28
$ws4.Range("F38").Value = 316
$ws4.Range("F42").Value = 278
$ws4.Range("F43").Value = 652
$ws4.Range("F45").Value = 46
